$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 285 (shifts ADL..PPT down by one row)
$ws.Rows.Item(285).Insert()

# Populate the newly inserted row 285 with the Austin, US entry
$ws.Cells.Item(285, 1).Value = "AUS"
$ws.Cells.Item(285, 2).Value = "Austin"
$ws.Cells.Item(285, 3).Value = 30.1975
$ws.Cells.Item(285, 4).Value = -97.6664
$ws.Cells.Item(285, 5).Value = "US"
$ws.Cells.Item(285, 6).Value = "North America"
$ws.Cells.Item(285, 7).Value = "Austin"

# Match the style of column A on the other data rows (bold, bordered, centered)
$ws.Cells.Item(286, 1).Copy()
$ws.Cells.Item(285, 1).PasteSpecial(-4122)
